{"js": "// Load all paragraphs in the body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// First paragraph's text becomes \"WORK\".\nparagraphs.items[0].insertText(\"WORK\", Word.InsertLocation.replace);\n\n// All the paragraphs in between (the \"Elementary/middle school\" block\n// through \"College:\") are removed, leaving only the trailing paragraph\n// (the one holding the _GoBack bookmark) after the first paragraph.\nfor (let i = paragraphs.items.length - 2; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# First paragraph's text becomes \"WORK\".\n$d.Paragraphs(1).Range.Text = \"WORK\"\n\n# Remove all the paragraphs between the first (\"WORK\") paragraph and the\n# last paragraph (the one carrying the _GoBack bookmark) -- this drops the\n# \"Elementary/middle school\" list block, \"High school: \", and \"College:\"\n# paragraphs entirely.\nfor ($i = $d.Paragraphs.Count - 1; $i -ge 2; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n"}
